# Update betting odds values on Sheet1 to reflect the latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Ind. Medellin - Llaneros) ---
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.9
$ws.Range("AR4").Value = 2.55
$ws.Range("AS4").Value = 1.49

# --- Row 17 (Difaa El Jadidi - Union Touarga) ---
$ws.Range("G17").Value = 2.1
$ws.Range("I17").Value = 3.4
$ws.Range("J17").Value = 3
$ws.Range("L17").Value = 4.33
$ws.Range("Z17").Value = 9.5
$ws.Range("AA17").Value = 9.5
$ws.Range("AB17").Value = 21
$ws.Range("AJ17").Value = 8.5
$ws.Range("AK17").Value = 17
$ws.Range("AL17").Value = 13
$ws.Range("AM17").Value = 41
$ws.Range("AN17").Value = 34

# --- Row 18 (AD Tarma - Sport Boys) ---
$ws.Range("G18").Value = 1.38
$ws.Range("I18").Value = 7.5
$ws.Range("J18").Value = 1.91
$ws.Range("K18").Value = 2.38
$ws.Range("L18").Value = 8
$ws.Range("O18").Value = 1.25
$ws.Range("P18").Value = 3.75
$ws.Range("Q18").Value = 1.88
$ws.Range("R18").Value = 1.98
$ws.Range("AA18").Value = 9
$ws.Range("AB18").Value = 8.5
$ws.Range("AD18").Value = 34
$ws.Range("AJ18").Value = 17
$ws.Range("AK18").Value = 41
$ws.Range("AL18").Value = 23
$ws.Range("AM18").Value = 101
